# Auto-generated Excel COM-interop script
# Applies updated FFXIV market-board derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 345.36365
$ws.Range("I19").Value = 309
$ws.Range("J19").Value = 375.66666
$ws.Range("K19").Value = 309
$ws.Range("L19").Value = 375.66666
$ws.Range("M19").Value = -134
$ws.Range("N19").Value = -725.66666
# Row 28
$ws.Range("H28").Value = 3373.25
$ws.Range("I28").Value = 929.7273
$ws.Range("J28").Value = 8749
$ws.Range("K28").Value = 929.7273
$ws.Range("L28").Value = 8749
$ws.Range("M28").Value = -444.7273
$ws.Range("N28").Value = -9719
# Row 40
$ws.Range("H40").Value = 4165.074
$ws.Range("I40").Value = 3274.7896
$ws.Range("K40").Value = 3274.7896
$ws.Range("M40").Value = -3099.7896
# Row 41
$ws.Range("H41").Value = 3145.3333
$ws.Range("I41").Value = 4500
$ws.Range("J41").Value = 2874.4
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 2874.4
$ws.Range("M41").Value = -4060
$ws.Range("N41").Value = -3754.4
# Row 106
$ws.Range("H106").Value = 5666.3335
$ws.Range("I106").Value = 5666.3335
$ws.Range("K106").Value = 5666.3335
$ws.Range("M106").Value = -5035.3335

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7091.269
$ws.Range("I32").Value = 4624.913
$ws.Range("K32").Value = 4624.913
$ws.Range("M32").Value = -4337.913
# Row 132
$ws.Range("H132").Value = 3369.6
$ws.Range("I132").Value = 3369.6
$ws.Range("K132").Value = 10108.8
$ws.Range("M132").Value = -7578.799999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1994.05
$ws.Range("I20").Value = 2047
$ws.Range("K20").Value = 2047
$ws.Range("M20").Value = -1800
# Row 86
$ws.Range("H86").Value = 3657.0527
$ws.Range("I86").Value = 2362.6428
$ws.Range("J86").Value = 7281.4
$ws.Range("K86").Value = 2362.6428
$ws.Range("L86").Value = 7281.4
$ws.Range("M86").Value = -1239.6428
$ws.Range("N86").Value = -9527.4
# Row 89
$ws.Range("H89").Value = 3657.0527
$ws.Range("I89").Value = 2362.6428
$ws.Range("J89").Value = 7281.4
$ws.Range("K89").Value = 11813.214
$ws.Range("L89").Value = 36407
$ws.Range("M89").Value = -6197.214
$ws.Range("N89").Value = -47639
# Row 107
$ws.Range("H107").Value = 3493.36
$ws.Range("I107").Value = 843.17645
$ws.Range("K107").Value = 843.17645
$ws.Range("M107").Value = 1076.82355
# Row 134
$ws.Range("H134").Value = 1729.6666
$ws.Range("I134").Value = 1729.6666
$ws.Range("K134").Value = 5188.9998
$ws.Range("M134").Value = -2653.9998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 4177.3477
$ws.Range("I132").Value = 3698.8462
$ws.Range("K132").Value = 11096.5386
$ws.Range("M132").Value = -8566.5386
# Row 134
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 2500
$ws.Range("K134").Value = 7500
$ws.Range("M134").Value = -4965

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 224.26315
$ws.Range("J12").Value = 223.64285
$ws.Range("L12").Value = 670.9285500000001
$ws.Range("N12").Value = -1016.92855
# Row 17
$ws.Range("H17").Value = 504.4
$ws.Range("I17").Value = 50.636364
$ws.Range("J17").Value = 1752.25
$ws.Range("K17").Value = 151.909092
$ws.Range("L17").Value = 5256.75
$ws.Range("M17").Value = 17.09090800000001
$ws.Range("N17").Value = -5594.75
# Row 98
$ws.Range("H98").Value = 161.88889
$ws.Range("I98").Value = 147.5
$ws.Range("J98").Value = 166
$ws.Range("K98").Value = 442.5
$ws.Range("L98").Value = 498
$ws.Range("M98").Value = 1055.5
$ws.Range("N98").Value = -3494

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 10041
$ws.Range("I46").Value = 10041
$ws.Range("K46").Value = 10041
$ws.Range("M46").Value = -9885
# Row 70
$ws.Range("H70").Value = 5855.5
$ws.Range("I70").Value = 5855.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5855.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5585.5
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 5855.5
$ws.Range("I73").Value = 5855.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5855.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4919.5
$ws.Range("N73").ClearContents()
# Row 80
$ws.Range("H80").Value = 4328.3335
$ws.Range("I80").Value = 4328.3335
$ws.Range("K80").Value = 4328.3335
$ws.Range("M80").Value = -3330.3335
# Row 83
$ws.Range("H83").Value = 4328.3335
$ws.Range("I83").Value = 4328.3335
$ws.Range("K83").Value = 21641.6675
$ws.Range("M83").Value = -16649.6675
# Row 97
$ws.Range("H97").Value = 802.4167
$ws.Range("I97").Value = 764.6667
$ws.Range("K97").Value = 764.6667
$ws.Range("M97").Value = -268.6667
# Row 107
$ws.Range("H107").Value = 250.2
$ws.Range("I107").Value = 170.4
$ws.Range("J107").Value = 330
$ws.Range("K107").Value = 170.4
$ws.Range("L107").Value = 330
$ws.Range("M107").Value = 1749.6
$ws.Range("N107").Value = -4170
# Row 132
$ws.Range("H132").Value = 5102.1113
$ws.Range("I132").Value = 4924.9287
$ws.Range("J132").Value = 5722.25
$ws.Range("K132").Value = 14774.7861
$ws.Range("L132").Value = 17166.75
$ws.Range("M132").Value = -12244.7861
$ws.Range("N132").Value = -22226.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9080
$ws.Range("I7").Value = 8244.5
$ws.Range("J7").Value = 9497.75
$ws.Range("K7").Value = 8244.5
$ws.Range("L7").Value = 9497.75
$ws.Range("M7").Value = -8132.5
$ws.Range("N7").Value = -9721.75
# Row 46
$ws.Range("H46").Value = 4115.72
$ws.Range("I46").Value = 3125.111
$ws.Range("J46").Value = 4672.9375
$ws.Range("K46").Value = 3125.111
$ws.Range("L46").Value = 4672.9375
$ws.Range("M46").Value = -2937.111
$ws.Range("N46").Value = -5048.9375
# Row 93
$ws.Range("H93").Value = 2010.8334
$ws.Range("I93").Value = 1992.2
$ws.Range("K93").Value = 1992.2
$ws.Range("M93").Value = -744.2
# Row 103
$ws.Range("H103").Value = 21820.2
$ws.Range("J103").Value = 21820.2
$ws.Range("L103").Value = 21820.2
$ws.Range("N103").Value = -24164.2
# Row 105
$ws.Range("H105").Value = 37666
$ws.Range("J105").Value = 37666
$ws.Range("L105").Value = 37666
$ws.Range("N105").Value = -44654
# Row 126
$ws.Range("H126").Value = 9080
$ws.Range("I126").Value = 8244.5
$ws.Range("J126").Value = 9497.75
$ws.Range("K126").Value = 24733.5
$ws.Range("L126").Value = 28493.25
$ws.Range("M126").Value = -22263.5
$ws.Range("N126").Value = -33433.25
# Row 132
$ws.Range("H132").Value = 2499.3333
$ws.Range("I132").Value = 2998.6667
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 8996.000100000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -6466.000100000001
$ws.Range("N132").Value = -11060
# Row 136
$ws.Range("H136").Value = 2916.6667
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 39089.715
$ws.Range("J54").Value = 54130.75
$ws.Range("L54").Value = 54130.75
$ws.Range("N54").Value = -55170.75
# Row 107
$ws.Range("H107").Value = 889.25
$ws.Range("I107").Value = 763.44446
$ws.Range("K107").Value = 2290.33338
$ws.Range("M107").Value = -370.33338
# Row 126
$ws.Range("H126").Value = 3243.4
$ws.Range("I126").Value = 1342.875
$ws.Range("K126").Value = 4028.625
$ws.Range("M126").Value = -1558.625

